$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Matriz_Resultados
# The 6h comparison between "Sieve Bootstrap" and "DeepAR" flips
# from significant (1 / -1) to not-significant (0) after the DM
# test correction, and DeepAR's row becomes all zeros.
# ---------------------------------------------------------------
$wsMatriz = $wb.Worksheets.Item("Matriz_Resultados")
$wsMatriz.Range("E2").Value = 0
$wsMatriz.Range("E3").Value = 0
$wsMatriz.Range("E4").Value = 0
$wsMatriz.Range("B5").Value = 0
$wsMatriz.Range("C5").Value = 0
$wsMatriz.Range("D5").Value = 0

# ---------------------------------------------------------------
# Sheet: P_valores
# Recomputed p-values for the corrected Diebold-Mariano / HLN test.
# ---------------------------------------------------------------
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.002885802965638673
$wsP.Range("D2").Value = 0.001514757432216962
$wsP.Range("E2").Value = 0.009871223014168651

$wsP.Range("B3").Value = 0.002885802965638673
$wsP.Range("D3").Value = 0.002351742784834299
$wsP.Range("E3").Value = 0.01374054425335625

$wsP.Range("B4").Value = 0.001514757432216962
$wsP.Range("C4").Value = 0.002351742784834299
$wsP.Range("E4").Value = 0.04198433918428757

$wsP.Range("B5").Value = 0.009871223014168651
$wsP.Range("C5").Value = 0.01374054425335625
$wsP.Range("D5").Value = 0.04198433918428757

# ---------------------------------------------------------------
# Sheet: Estadisticos_HLN_DM
# Recomputed HLN-corrected DM statistics.
# ---------------------------------------------------------------
$wsE = $wb.Worksheets.Item("Estadisticos_HLN_DM")

$wsE.Range("C2").Value = -3.351582633019906
$wsE.Range("D2").Value = -3.620494141815038
$wsE.Range("E2").Value = -2.82445035501861

$wsE.Range("B3").Value = 3.351582633019906
$wsE.Range("D3").Value = -3.437363438463411
$wsE.Range("E3").Value = -2.677967864749372

$wsE.Range("B4").Value = 3.620494141815038
$wsE.Range("C4").Value = 3.437363438463411
$wsE.Range("E4").Value = -2.159439046062592

$wsE.Range("B5").Value = 2.82445035501861
$wsE.Range("C5").Value = 2.677967864749372
$wsE.Range("D5").Value = 2.159439046062592

# ---------------------------------------------------------------
# Sheet: Resumen_Modelos
# Updated win/loss/tie counts and win-rate strings after the fix.
# Percent-looking cells must stay as literal text, not be
# reinterpreted by Excel as numeric percentages, so force the
# cell to Text format before writing the string.
# ---------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Resumen_Modelos")

$wsR.Range("B2").Value = 2
$wsR.Range("D2").Value = 1
$wsR.Range("E2").NumberFormat = "@"
$wsR.Range("E2").Value = "66.7%"

$wsR.Range("B3").Value = 1
$wsR.Range("D3").Value = 1
$wsR.Range("E3").NumberFormat = "@"
$wsR.Range("E3").Value = "33.3%"

$wsR.Range("B4").Value = 0
$wsR.Range("D4").Value = 1
$wsR.Range("E4").NumberFormat = "@"
$wsR.Range("E4").Value = "0.0%"

$wsR.Range("C5").Value = 0
$wsR.Range("D5").Value = 3
